# Auto-generated edit script applying numeric corrections to H:N columns
# across multiple sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ALC")
$ws2 = $wb.Worksheets.Item("ARM")
$ws3 = $wb.Worksheets.Item("BSM")
$ws4 = $wb.Worksheets.Item("CRP")
$ws5 = $wb.Worksheets.Item("CUL")
$ws6 = $wb.Worksheets.Item("GSM")
$ws7 = $wb.Worksheets.Item("LTW")
$ws8 = $wb.Worksheets.Item("WVR")

# @@ -2506,22 +2506,22 @@  (sheet ALC)
$ws1.Range("H38").Value = 14478.429
$ws1.Range("I38").Value = 224.83333
$ws1.Range("K38").Value = 674.49999
$ws1.Range("M38").Value = -302.49999

# @@ -3703,22 +3703,22 @@  (sheet ALC)
$ws1.Range("H62").Value = 8075
$ws1.Range("I62").Value = 7719.1665
$ws1.Range("K62").Value = 7719.1665
$ws1.Range("M62").Value = -7095.1665

# @@ -3853,22 +3853,22 @@  (sheet ALC)
$ws1.Range("H65").Value = 8075
$ws1.Range("I65").Value = 7719.1665
$ws1.Range("K65").Value = 38595.8325
$ws1.Range("M65").Value = -35475.8325

# @@ -5221,25 +5221,25 @@  (sheet ALC)
$ws1.Range("H92").Value = 16667198
$ws1.Range("I92").Value = 18519048
$ws1.Range("J92").Value = 556.6667
$ws1.Range("K92").Value = 18519048
$ws1.Range("L92").Value = 556.6667
$ws1.Range("M92").Value = -18517800
$ws1.Range("N92").Value = -3052.6667

# @@ -5423,25 +5423,25 @@  (sheet ALC)
$ws1.Range("H96").Value = 1350.1428
$ws1.Range("I96").Value = 408.66666
$ws1.Range("J96").Value = 2056.25
$ws1.Range("K96").Value = 1225.99998
$ws1.Range("L96").Value = 6168.75
$ws1.Range("M96").Value = 147.0000199999999
$ws1.Range("N96").Value = -8914.75

# @@ -5934,22 +5934,22 @@  (sheet ALC)
$ws1.Range("H106").Value = 7411536.5
$ws1.Range("I106").Value = 8337766
$ws1.Range("K106").Value = 8337766
$ws1.Range("M106").Value = -8337135

# @@ -6390,22 +6390,22 @@  (sheet ALC)
$ws1.Range("H115").Value = 552.7143
$ws1.Range("I115").Value = 561.6667
$ws1.Range("K115").Value = 1685.0001
$ws1.Range("M115").Value = -118.0001

# @@ -7180,25 +7180,25 @@  (sheet ALC)
$ws1.Range("H131").Value = 2359.2
$ws1.Range("I131").Value = 1580.5
$ws1.Range("J131").Value = 5474
$ws1.Range("K131").Value = 4741.5
$ws1.Range("L131").Value = 16422
$ws1.Range("M131").Value = 298.5
$ws1.Range("N131").Value = -26502

# @@ -7232,22 +7232,22 @@  (sheet ALC)
$ws1.Range("H132").Value = 2421.7058
$ws1.Range("I132").Value = 2260.875
$ws1.Range("K132").Value = 6782.625
$ws1.Range("M132").Value = -4252.625

# @@ -7474,25 +7474,25 @@  (sheet ALC)
$ws1.Range("H137").Value = 2969716
$ws1.Range("I137").Value = 70225.164
$ws1.Range("J137").Value = 7940272
$ws1.Range("K137").Value = 210675.492
$ws1.Range("L137").Value = 23820816
$ws1.Range("M137").Value = -208125.492
$ws1.Range("N137").Value = -23825916

# @@ -7526,25 +7526,22 @@  (sheet ALC)
$ws1.Range("H138").Value = 2929.67
$ws1.Range("I138").Value = 0
$ws1.Range("J138").Value = 2929.67
$ws1.Range("K138").Value = 0
$ws1.Range("L138").Value = 8789.01
$ws1.Range("M138").ClearContents()
$ws1.Range("N138").Value = -19069.01

# @@ -7673,7 +7670,7 @@  (sheet ALC)
$ws1.Range("H141").Value = 5996.75

# @@ -7969,22 +7966,25 @@  (sheet ARM)
$ws2.Range("H5").Value = 309.8
$ws2.Range("I5").Value = 333.1111
$ws2.Range("J5").Value = 100
$ws2.Range("K5").Value = 333.1111
$ws2.Range("L5").Value = 100
$ws2.Range("M5").Value = -221.1111
$ws2.Range("N5").Value = -324

# @@ -9932,22 +9932,22 @@  (sheet ARM)
$ws2.Range("H45").Value = 4449.6895
$ws2.Range("I45").Value = 4821.8887
$ws2.Range("K45").Value = 4821.8887
$ws2.Range("M45").Value = -4444.8887

# @@ -10805,25 +10805,25 @@  (sheet ARM)
$ws2.Range("H63").Value = 4500
$ws2.Range("J63").Value = 5500
$ws2.Range("L63").Value = 5500
$ws2.Range("N63").Value = -6872

# @@ -10952,25 +10952,25 @@  (sheet ARM)
$ws2.Range("H66").Value = 4500
$ws2.Range("J66").Value = 5500
$ws2.Range("L66").Value = 27500
$ws2.Range("N66").Value = -34364

# @@ -11338,22 +11338,22 @@  (sheet ARM)
$ws2.Range("H74").Value = 2165.484
$ws2.Range("I74").Value = 2154.8076
$ws2.Range("K74").Value = 2154.8076
$ws2.Range("M74").Value = -1280.8076

# @@ -11482,22 +11482,22 @@  (sheet ARM)
$ws2.Range("H77").Value = 2165.484
$ws2.Range("I77").Value = 2154.8076
$ws2.Range("K77").Value = 10774.038
$ws2.Range("M77").Value = -6406.038

# @@ -12704,25 +12704,25 @@  (sheet ARM)
$ws2.Range("H102").Value = 1804.8
$ws2.Range("I102").Value = 1539.3334
$ws2.Range("J102").Value = 2866.6667
$ws2.Range("K102").Value = 1539.3334
$ws2.Range("L102").Value = 2866.6667
$ws2.Range("M102").Value = 82.66660000000002
$ws2.Range("N102").Value = -6110.6667

# @@ -14865,22 +14865,25 @@  (sheet BSM)
$ws3.Range("H4").Value = 309.8
$ws3.Range("I4").Value = 333.1111
$ws3.Range("J4").Value = 100
$ws3.Range("K4").Value = 333.1111
$ws3.Range("L4").Value = 100
$ws3.Range("M4").Value = -218.1111
$ws3.Range("N4").Value = -330

# @@ -16611,22 +16614,22 @@  (sheet BSM)
$ws3.Range("H40").Value = 56448
$ws3.Range("J40").Value = 56448
$ws3.Range("L40").Value = 56448
$ws3.Range("N40").Value = -56978

# @@ -19775,22 +19778,22 @@  (sheet BSM)
$ws3.Range("H105").Value = 2081.158
$ws3.Range("I105").Value = 1782.9333
$ws3.Range("K105").Value = 1782.9333
$ws3.Range("M105").Value = -35.93329999999992

# @@ -22341,22 +22344,22 @@  (sheet CRP)
$ws4.Range("H16").Value = 1229.4286
$ws4.Range("I16").Value = 1383.7273
$ws4.Range("K16").Value = 1383.7273
$ws4.Range("M16").Value = -1096.7273

# @@ -23088,25 +23091,25 @@  (sheet CRP)
$ws4.Range("H31").Value = 5231.6
$ws4.Range("I31").Value = 3301
$ws4.Range("J31").Value = 5803.6294
$ws4.Range("K31").Value = 3301
$ws4.Range("L31").Value = 5803.6294
$ws4.Range("M31").Value = -3006
$ws4.Range("N31").Value = -6393.6294

# @@ -23238,25 +23241,25 @@  (sheet CRP)
$ws4.Range("H34").Value = 5231.6
$ws4.Range("I34").Value = 3301
$ws4.Range("J34").Value = 5803.6294
$ws4.Range("K34").Value = 3301
$ws4.Range("L34").Value = 5803.6294
$ws4.Range("M34").Value = -3099
$ws4.Range("N34").Value = -6207.6294

# @@ -26693,22 +26696,22 @@  (sheet CRP)
$ws4.Range("H105").Value = 1823.4117
$ws4.Range("I105").Value = 1464.1428
$ws4.Range("K105").Value = 1464.1428
$ws4.Range("M105").Value = 282.8571999999999

# @@ -27091,22 +27094,22 @@  (sheet CRP)
$ws4.Range("H113").Value = 1229.4286
$ws4.Range("I113").Value = 1383.7273
$ws4.Range("K113").Value = 1383.7273
$ws4.Range("M113").Value = 786.2727

# @@ -27529,25 +27532,25 @@  (sheet CRP)
$ws4.Range("H122").Value = 5375
$ws4.Range("I122").Value = 4298
$ws4.Range("J122").Value = 7067.4287
$ws4.Range("K122").Value = 12894
$ws4.Range("L122").Value = 21202.2861
$ws4.Range("M122").Value = -10444
$ws4.Range("N122").Value = -26102.2861

# @@ -28451,22 +28454,22 @@  (sheet CRP)
$ws4.Range("H141").Value = 453327.16
$ws4.Range("J141").Value = 453327.16
$ws4.Range("L141").Value = 453327.16
$ws4.Range("N141").Value = -463687.16

# @@ -34567,25 +34570,25 @@  (sheet CUL)
$ws5.Range("H122").Value = 573.375
$ws5.Range("I122").Value = 472.83334
$ws5.Range("J122").Value = 633.7
$ws5.Range("K122").Value = 4255.50006
$ws5.Range("L122").Value = 5703.3
$ws5.Range("M122").Value = -1805.50006
$ws5.Range("N122").Value = -10603.3

# @@ -41846,22 +41849,22 @@  (sheet GSM)
$ws6.Range("H129").Value = 119995
$ws6.Range("J129").Value = 119995
$ws6.Range("L129").Value = 119995
$ws6.Range("N129").Value = -129995

# @@ -46957,22 +46960,22 @@  (sheet LTW)
$ws7.Range("H92").Value = 40379
$ws7.Range("J92").Value = 40379
$ws7.Range("L92").Value = 40379
$ws7.Range("N92").Value = -45371

# @@ -48418,22 +48421,22 @@  (sheet LTW)
$ws7.Range("H122").Value = 10094.5
$ws7.Range("I122").Value = 9197.75
$ws7.Range("K122").Value = 27593.25
$ws7.Range("M122").Value = -25143.25

# @@ -52405,25 +52408,25 @@  (sheet WVR)
$ws8.Range("H62").Value = 5037.5
$ws8.Range("I62").Value = 3166.6667
$ws8.Range("J62").Value = 5661.1113
$ws8.Range("K62").Value = 3166.6667
$ws8.Range("L62").Value = 5661.1113
$ws8.Range("M62").Value = -2542.6667
$ws8.Range("N62").Value = -6909.1113

# @@ -52555,25 +52558,25 @@  (sheet WVR)
$ws8.Range("H65").Value = 5037.5
$ws8.Range("I65").Value = 3166.6667
$ws8.Range("J65").Value = 5661.1113
$ws8.Range("K65").Value = 15833.3335
$ws8.Range("L65").Value = 28305.5565
$ws8.Range("M65").Value = -12713.3335
$ws8.Range("N65").Value = -34545.5565

# @@ -54883,22 +54886,22 @@  (sheet WVR)
$ws8.Range("H113").Value = 361.22223
$ws8.Range("I113").Value = 310.2
$ws8.Range("K113").Value = 930.5999999999999
$ws8.Range("M113").Value = 1239.4

# @@ -56016,22 +56019,22 @@  (sheet WVR)
$ws8.Range("H136").Value = 2907.5386
$ws8.Range("I136").Value = 1981.7273
$ws8.Range("K136").Value = 5945.1819
$ws8.Range("M136").Value = -3395.1819
